# navi bar for payment handling
#
# Updates the "Current Iteration" navigation bar on slide 2 (iteration
# number + date ranges) and refreshes the fixed "date last edited" text
# that is stamped into the slide master and every slide layout's Date
# Placeholder (9/3/2018 -> 9/17/2018).

$p = $ppt.ActivePresentation

# PowerPoint's text-diffing will happily keep a long common sub-string
# between the old and new text alive as a separate run (e.g. turning a
# single run into "17" + " Sep 2018 " runs) which does not match how the
# source file was actually edited (a single run holding the whole
# paragraph). Routing every edit through an unrelated placeholder value
# first guarantees there is no overlap with the previous text, so the
# final assignment always collapses back down to one clean run.
# NOTE: this interpreter's functions only bind positional parameters
# reliably, so always call Set-CleanText with positional args.
function Set-CleanText {
    param($TextRange, $NewText)

    $TextRange.Text = "#_TMP_PLACEHOLDER_TEXT_#"
    $TextRange.Text = $NewText
}

# ---------------------------------------------------------------------
# 1. Slide 2 "Current Iteration" navigation bar
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Iteration number box ("TextBox 22"): 6 -> 7
$iterBox = $s2.Shapes.Item(5)
Set-CleanText $iterBox.TextFrame.TextRange "7"

# Iteration date range box ("TextBox 24"):
#   "2 Sep 2018 " / "   -  16 Sep 2018"  ->
#   "17 Sep 2018 " / "   -  30 Sep 2018"
$dateRangeBox = $s2.Shapes.Item(7)
$dash = [char]0x2013
Set-CleanText $dateRangeBox.TextFrame.TextRange.Paragraphs(1,1) "17 Sep 2018 "
Set-CleanText $dateRangeBox.TextFrame.TextRange.Paragraphs(2,1) ("   " + $dash + "  30 Sep 2018")

# The three meeting-schedule date boxes further down slide 2
Set-CleanText $s2.Shapes.Item(36).TextFrame.TextRange "21 Sep 2018, Fri"
Set-CleanText $s2.Shapes.Item(37).TextFrame.TextRange "22 Sep 2018, Sat"
Set-CleanText $s2.Shapes.Item(38).TextFrame.TextRange "21 Sep 2018, Fri"

# ---------------------------------------------------------------------
# 2. Fixed date stamp (9/3/2018 -> 9/17/2018) on the slide master and
#    every slide layout's Date Placeholder.
# ---------------------------------------------------------------------
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        Set-CleanText $sh.TextFrame.TextRange "9/17/2018"
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $sh = $layout.Shapes.Item($si)
        if ($sh.Name -like "Date Placeholder*") {
            Set-CleanText $sh.TextFrame.TextRange "9/17/2018"
        }
    }
}
